$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 11, shifting the existing row 11 (and below) down.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with data (same content as row 10).
$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "dio"
$ws.Range("C11").Value = "world"
$ws.Range("D11").Value = "kk"
$ws.Range("E11").Value = 1998
$ws.Range("F11").Value = 888
$ws.Range("G11").Value = 55

# Match the final selection state: entire row 11 selected, active cell A11.
$ws.Rows.Item(11).Select()
